$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reword the graduation-count and graduation-rate column headers to spell
# out the numbers instead of using digits.
$ws.Range("C1").Value = "four year graduation count"
$ws.Range("D1").Value = "six year graduation count"
$ws.Range("E1").Value = "eight year graduation count"
$ws.Range("F1").Value = "four year graduation rate"
$ws.Range("G1").Value = "six year graduation rate"
$ws.Range("H1").Value = "eight year graduation rate"

# The header row grows taller (likely so the new, longer header text wraps).
$ws.Rows.Item(1).RowHeight = 30

# Move the active selection to D9, matching where the author was last
# working.
$ws.Range("D9").Select()
